# Regenerate the "K" (strikeouts) column (column G) using correct values
# instead of the old "Strike#" derived values, for each outing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(3,1,4,2,5,5,3,3,0,4,0,7,8,4,4,12,6,2,2,5,6,3,3,5,5,3,8,4,4,3,4,6,4)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
